# Refresh the crypto price/ranking snapshot (GitHub Actions symbol-list update).
# The source feed re-ranked a handful of coins (rows 18-24 shift by one slot as
# "One" jumps up to #17) and refreshed several Price / Volume(1h) column values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "242.70" },
    @{ Cell = "D3";  Value = "23.12" },
    @{ Cell = "D4";  Value = "5.422" },
    @{ Cell = "D5";  Value = "0.05923" },
    @{ Cell = "D6";  Value = "3.451" },
    @{ Cell = "D7";  Value = "6.553" },
    @{ Cell = "D8";  Value = "0.8131" },
    @{ Cell = "D9";  Value = "0.9068" },
    @{ Cell = "D10"; Value = "0.1409" },
    @{ Cell = "D11"; Value = "0.07452" },
    @{ Cell = "D12"; Value = "0.03300" },
    @{ Cell = "D13"; Value = "0.03057" },
    @{ Cell = "D14"; Value = "0.09340" },
    @{ Cell = "D15"; Value = "3.850" },
    @{ Cell = "D16"; Value = "0.001589" },
    @{ Cell = "D17"; Value = "0.04675" },

    # Row 18: TigerCash -> One
    @{ Cell = "B18"; Value = "One" },
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one" },
    @{ Cell = "D18"; Value = "0.0005939" },
    @{ Cell = "E18"; Value = "17OneONE" },

    # Row 19: HotbitToken -> TigerCash
    @{ Cell = "B19"; Value = "TigerCash" },
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" },
    @{ Cell = "D19"; Value = "0.006081" },
    @{ Cell = "E19"; Value = "18TigerCashTCH" },

    # Row 20: BitKan -> HotbitToken
    @{ Cell = "B20"; Value = "HotbitToken" },
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb" },
    @{ Cell = "D20"; Value = "0.004982" },
    @{ Cell = "E20"; Value = "19HotbitTokenHTB" },

    # Row 21: NitroEx -> BitKan
    @{ Cell = "B21"; Value = "BitKan" },
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan" },
    @{ Cell = "D21"; Value = "0.0009864" },
    @{ Cell = "E21"; Value = "20BitKanKAN" },

    # Row 22: LEO -> NitroEx
    @{ Cell = "B22"; Value = "NitroEx" },
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx" },
    @{ Cell = "D22"; Value = "0.00008998" },
    @{ Cell = "E22"; Value = "21NitroExNTX" },

    # Row 23: BTSEToken -> LEO
    @{ Cell = "B23"; Value = "LEO" },
    @{ Cell = "C23"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" },
    @{ Cell = "D23"; Value = "3.604" },
    @{ Cell = "E23"; Value = "22LEOLEO" },

    # Row 24: One -> BTSEToken
    @{ Cell = "B24"; Value = "BTSEToken" },
    @{ Cell = "C24"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse" },
    @{ Cell = "D24"; Value = "2.135" },
    @{ Cell = "E24"; Value = "23BTSETokenBTSE" },

    @{ Cell = "D25"; Value = "0.3240" },
    @{ Cell = "D27"; Value = "0.0002899" },
    @{ Cell = "D40"; Value = "0.04015" },
    @{ Cell = "D41"; Value = "0.006213" },
    @{ Cell = "E41"; Value = "40KickTokenKICKBestin24h" },
    @{ Cell = "D43"; Value = "0.002999" },
    @{ Cell = "D44"; Value = "0.008095" },
    @{ Cell = "D45"; Value = "0.00005244" },
    @{ Cell = "D48"; Value = "0.8156" },
    @{ Cell = "D49"; Value = "0.002254" },
    @{ Cell = "D50"; Value = "0.00002100" },
    @{ Cell = "D51"; Value = "0.0002000" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (leading/trailing zeros,
    # very small decimals) round-trip verbatim instead of being normalized
    # into a Double (which would drop trailing zeros / use sci notation).
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
